$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = "974 BRAHMS CT"
$ws.Range("D5").Value = "975 BRAHMS CT"
$ws.Range("D6").Value = "976 BRAHMS CT"
